# Update significance-marked p-value strings to remove the stray space
# between the numeric value and the trailing asterisk(s).

$wb = $excel.ActiveWorkbook

$rich = $wb.Worksheets.Item("rich")
$rich.Range("F4").Value = "0.001**"
$rich.Range("F9").Value = "0.008**"
$rich.Range("F11").Value = "0.01*"
$rich.Range("F12").Value = "0.044*"
$rich.Range("F14").Value = "0.002**"

$even = $wb.Worksheets.Item("even")
$even.Range("F2").Value = "<0.001***"
$even.Range("F4").Value = "<0.001***"
$even.Range("F9").Value = "0.005**"
$even.Range("F10").Value = "0.008**"
$even.Range("F11").Value = "0.01*"
$even.Range("F12").Value = "0.001**"

$invSim = $wb.Worksheets.Item("invSim")
$invSim.Range("F2").Value = "0.015*"
$invSim.Range("F4").Value = "<0.001***"
$invSim.Range("F13").Value = "0.041*"
$invSim.Range("F15").Value = "0.003**"
